$wb = $excel.ActiveWorkbook

# Each worksheet gets one new data row appended at the bottom, matching the
# same column layout as the existing rows (A: timestamp, B-E: text fields,
# F/H/I: numbers, G: usually a big number but stored as text on sheet 2).

# --- Sheet 1: ROW50-FE-LIFTER -> new row 59 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r1 = 59
$ws1.Range("A$r1").Value = 45753.71616379629
$ws1.Range("A$r1").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B$r1").Value = "0x01,0x90"
$ws1.Range("C$r1").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Range("D$r1").Value = "0x01,0x5a"
$ws1.Range("E$r1").Value = "0xe"
$ws1.Range("F$r1").Value = 400
$ws1.Range("G$r1").Value = [double]"5.68631262647114e+23"
$ws1.Range("H$r1").Value = 346
$ws1.Range("I$r1").Value = 14

# --- Sheet 2: ROW50-MID-LIFTER -> new row 61 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r2 = 61
$ws2.Range("A$r2").Value = 45753.68271990741
$ws2.Range("A$r2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B$r2").Value = "0x01,0x90 "
$ws2.Range("C$r2").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D$r2").Value = "0x01,0x62"
$ws2.Range("E$r2").Value = "0x19"
$ws2.Range("F$r2").Value = 400
# G61 holds the value as plain text (not numeric) to match the source row.
# Writing the digit string directly would auto-convert to a number, so build
# it as a text formula first and flatten it to a literal value via paste-values
# (keeps the default/general cell style, unlike forcing a "@" number format).
$ws2.Range("G$r2").Formula = "=""568631262647113771663628"""
$ws2.Range("G$r2").Copy() | Out-Null
$ws2.Range("G$r2").PasteSpecial(-4163) | Out-Null
$ws2.Range("H$r2").Value = 354
$ws2.Range("I$r2").Value = 25

# --- Sheet 3: ROW11-FE-LIFTER -> new row 59 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r3 = 59
$ws3.Range("A$r3").Value = 45753.74834099537
$ws3.Range("A$r3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B$r3").Value = "0x01,0x90"
$ws3.Range("C$r3").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D$r3").Value = "0x01,0x5a"
$ws3.Range("E$r3").Value = "0x14"
$ws3.Range("F$r3").Value = 400
$ws3.Range("G$r3").Value = [double]"5.68631262647114e+23"
$ws3.Range("H$r3").Value = 346
$ws3.Range("I$r3").Value = 20

# --- Sheet 4: ROW11-MID-LIFTER -> new row 59 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r4 = 59
$ws4.Range("A$r4").Value = 45753.88044846065
$ws4.Range("A$r4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B$r4").Value = "0x01,0x90"
$ws4.Range("C$r4").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D$r4").Value = "0x01,0x62"
$ws4.Range("E$r4").Value = "0x19"
$ws4.Range("F$r4").Value = 400
$ws4.Range("G$r4").Value = [double]"5.68631262647114e+23"
$ws4.Range("H$r4").Value = 354
$ws4.Range("I$r4").Value = 25
